$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # RUNMANAGER
$ws2 = $wb.Worksheets.Item(2)   # DATA

# ---------------------------------------------------------------------------
# Sheet "DATA" (sheet2): new column G ("menutext"), new rows 7 & 8
# (order of first-use chosen so new shared strings land at the same indices
#  as the target: menutext=25, Laptops=26, amazonTest=27, <description>=28)
# ---------------------------------------------------------------------------

# New column G header
$ws2.Cells.Item(1, 7).Value = "menutext"
$ws2.Cells.Item(1, 7).Font.Size = 22

# Fill existing rows 2-6 in the new column G with an empty, quote-prefixed cell
$ws2.Cells.Item(2, 7).Value = "'"
$ws2.Cells.Item(3, 7).Value = "'"
$ws2.Cells.Item(4, 7).Value = "'"
$ws2.Cells.Item(5, 7).Value = "'"
$ws2.Cells.Item(6, 7).Value = "'"

# New row 7: amazonTest / chrome / Laptops
$ws2.Cells.Item(7, 7).Value = "Laptops"
$ws2.Cells.Item(7, 1).Value = "amazonTest"
$ws2.Cells.Item(7, 2).Value = "yes"
$ws2.Cells.Item(7, 3).Value = "chrome"
$ws2.Cells.Item(7, 4).Value = "'"
$ws2.Cells.Item(7, 5).Value = "'"
$ws2.Cells.Item(7, 6).Value = "'"
$ws2.Range("A7:C7").Font.Size = 22
$ws2.Cells.Item(7, 7).Font.Size = 22
$ws2.Rows.Item(7).RowHeight = 28.8

# New row 8: amazonTest / firefox / Laptops
$ws2.Cells.Item(8, 1).Value = "amazonTest"
$ws2.Cells.Item(8, 2).Value = "yes"
$ws2.Cells.Item(8, 3).Value = "firefox"
$ws2.Cells.Item(8, 4).Value = "'"
$ws2.Cells.Item(8, 5).Value = "'"
$ws2.Cells.Item(8, 6).Value = "'"
$ws2.Cells.Item(8, 7).Value = "Laptops"
$ws2.Range("A8:C8").Font.Size = 22
$ws2.Cells.Item(8, 7).Font.Size = 22
$ws2.Rows.Item(8).RowHeight = 28.8

# Column widths: col A widened (best-fit-like), new col G sized
$ws2.Columns.Item(1).ColumnWidth = 31.8
$ws2.Columns.Item(7).ColumnWidth = 16.8

# ---------------------------------------------------------------------------
# Sheet "RUNMANAGER" (sheet1): execute flags flipped to "no", new test row 4
# ---------------------------------------------------------------------------

$ws1.Cells.Item(2, 3).Value = "no"
$ws1.Cells.Item(3, 3).Value = "no"

$ws1.Cells.Item(4, 1).Value = "amazonTest"
$ws1.Cells.Item(4, 2).Value = "To check whether amazon test is executed "
$ws1.Cells.Item(4, 3).Value = "yes"
$ws1.Cells.Item(4, 4).Value = "'1"
$ws1.Cells.Item(4, 5).Value = "'1"
$ws1.Range("A4:E4").Font.Size = 22
$ws1.Rows.Item(4).RowHeight = 28.8

# ---------------------------------------------------------------------------
# Page setup (sheet2 gains an explicit portrait print orientation)
# ---------------------------------------------------------------------------
$ws2.PageSetup.Orientation = 1

# ---------------------------------------------------------------------------
# Selections, matching the final cursor position on each sheet
# ---------------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("A5:XFD8").Select()

$ws2.Activate()
$ws2.Range("F5").Select()
